# Insert two new price-report rows (Vega Central Mapocho de Santiago - Cebollín)
# right after the existing row 1244, shifting all subsequent rows down by two.
# This mirrors the diff: dimension grows from A1:R1302 to A1:R1304, with a new
# "Primera" record and a new "Segunda" record (both dated 44939 / 2023-01-13)
# inserted before what used to be row 1245.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 1245 (pushes old rows 1245.. down to 1247..)
$ws.Rows.Item(1245).Insert()
$ws.Rows.Item(1245).Insert()

# Common (constant) column values shared by every data row in this sheet
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$catId     = 100112037
$categoria = "Cebollín"
$variedad  = "Sin especificar"
$unidad    = "`$/paquete 36 unidades"
$origen    = "Región Metropolitana"
$kgUnid    = 36
$clasif    = "Hortaliza"

function Set-CebollinRow {
    param($Row, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg)

    $ws.Cells.Item($Row, 1).Value  = $mercadoId
    $ws.Cells.Item($Row, 2).Value  = $mercado
    $ws.Cells.Item($Row, 3).Value  = $region
    $ws.Cells.Item($Row, 4).Value  = $Fecha
    $ws.Cells.Item($Row, 5).Value  = $codreg
    $ws.Cells.Item($Row, 6).Value  = $catId
    $ws.Cells.Item($Row, 7).Value  = $categoria
    $ws.Cells.Item($Row, 8).Value  = $variedad
    $ws.Cells.Item($Row, 9).Value  = $Calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 14).Value = $unidad
    $ws.Cells.Item($Row, 15).Value = $origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = $kgUnid
    $ws.Cells.Item($Row, 18).Value = $clasif
}

# New row 1245: Primera, 13-01-2023 (serial 44939)
Set-CebollinRow 1245 44939 "Primera" 430 3000 3000 3000 83

# New row 1246: Segunda, 13-01-2023 (serial 44939)
Set-CebollinRow 1246 44939 "Segunda" 340 2500 2500 2500 69
